$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column C (Förändrad) changes from 46066 to 46070 for all data rows (2..56)
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 3).Value = 46070
}

# 2. Rows 43 and 44 swap their content (A, B, G columns)
$ws.Range("A43").Value = "A 29982-2025"
$ws.Range("B43").Value = 45826
$ws.Range("G43").Value = 1.5

$ws.Range("A44").Value = "A 12874-2021"
$ws.Range("B44").Value = 44271.42787037037
$ws.Range("G44").Value = 6

# 3. Rows 46-52 get rotated: the last two entries (60731-2025 / 60733-2025)
#    move to the top (rows 46-47), and the remaining five entries shift down
#    keeping their relative order (rows 48-52). Row 45 is unaffected.
$ws.Range("A46").Value = "A 60731-2025"
$ws.Range("B46").Value = 45995
$ws.Range("G46").Value = 0.8

$ws.Range("A47").Value = "A 60733-2025"
$ws.Range("B47").Value = 45995
$ws.Range("G47").Value = 1.7
$ws.Range("F47").Value = ""

$ws.Range("A48").Value = "A 33246-2021"
$ws.Range("B48").Value = 44377
$ws.Range("G48").Value = 3.4
$ws.Range("F48").Value = ""

$ws.Range("A49").Value = "A 46919-2023"
$ws.Range("B49").Value = 45201
$ws.Range("G49").Value = 1.3
$ws.Range("F49").Value = "Kommuner"

$ws.Range("A50").Value = "A 504-2023"
$ws.Range("B50").Value = 44930
$ws.Range("G50").Value = 2
$ws.Range("F50").Value = ""

$ws.Range("A51").Value = "A 58592-2023"
$ws.Range("B51").Value = 45251
$ws.Range("G51").Value = 2.1

$ws.Range("A52").Value = "A 61336-2022"
$ws.Range("B52").Value = 44915
$ws.Range("G52").Value = 5.6
$ws.Range("F52").Value = "Kommuner"
